$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: D1 becomes "Department " (previously this label lived in F1; program/year
# columns are dropped entirely)
$ws.Range("D1").Value = "Department "

# Department values for each student row (rows 2-21), replacing the old
# program/year columns with a single Department column
$depts = @("CSE","CE","CE","ME","EE","EE","EE","EE","CSE","CSE","CE","ME","ME","EE","EE","CE","CE","CE","CE","CE")

for ($i = 0; $i -lt $depts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $depts[$i]
}

# The old E (program) and F (year) columns are removed entirely
$ws.Range("E1:F21").Clear()

# Slightly widen column D to fit the new header/values
$ws.Columns.Item(4).ColumnWidth = 10.44140625

# Match the saved selection state
[void]$ws.Range("E20").Select()
